$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.685.62'
$ws.Range("E2").Value = '  -1.82%  '
$ws.Range("D3").Value = '2.441.29'
$ws.Range("E3").Value = '  -1.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.00'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.70%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.530'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.17%  '
$ws.Range("D9").Value = '2.438.85'
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.97%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.350'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.05%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.67%  '
$ws.Range("D16").Value = '2.893.22'
$ws.Range("E16").Value = '  -1.52%  '
$ws.Range("D17").Value = '62.652.44'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Value = '2.441.60'
$ws.Range("E18").Value = '  -1.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.68%  '
$ws.Range("E20").Value = '  -3.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '330.62'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.10%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.81'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '633.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.01'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.98%  '
$ws.Range("D28").Value = '0.0₃0964'
$ws.Range("E28").Value = '  -8.19%  '
$ws.Range("D29").Value = '2.563.67'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.07'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.49%  '
$ws.Range("E33").Value = '  -2.03%  '
$ws.Range("E34").Value = '  -3.74%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.55%  '
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("E38").Value = '  -2.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.87%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.29'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.75%  '
$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.98'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.02%  '
$ws.Range("E42").Value = '  -4.89%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.53'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.37%  '
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.89%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.71'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0524'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.72%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.599'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.73'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.71%  '
$ws.Range("D51").Value = '0.0₆0236'
$ws.Range("E51").Value = '  +7.95%  '
